# Apply the weekly Fruta/Hortaliza update: the rows in the data table are
# re-shuffled (dates, quality, volume and price columns move between rows)
# while columns A,B,C,E,F,G,H,I,J,K,Q,R,T stay constant for every row.
#
# Mapping is: new row <- old row (source of D,L,M,N,O,P,S values)
#   2<-12  3<-2   4<-5   5<-13  6<-9   7<-10  8<-6
#   9<-7   10<-8  11<-14 12<-3  13<-4  14<-11 15<-17
#   16<-15 17<-16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary row-to-row.
$cols = @(4, 12, 13, 14, 15, 16, 19)   # D, L, M, N, O, P, S

# Snapshot the "before" values for every row/column we need, so that
# writes to earlier rows don't clobber data still needed for later rows.
$snapshot = @{}
for ($r = 2; $r -le 17; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# new row -> source (old) row
$mapping = @{
    2  = 12
    3  = 2
    4  = 5
    5  = 13
    6  = 9
    7  = 10
    8  = 6
    9  = 7
    10 = 8
    11 = 14
    12 = 3
    13 = 4
    14 = 11
    15 = 17
    16 = 15
    17 = 16
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
